$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Description" note tag should no longer carry the type="local" attribute.
# The separators in the original text are non-breaking spaces (U+00A0), so we
# rebuild the string explicitly to preserve that formatting.
$nbsp = [char]0x00A0
$newNoteTag = "<mods:note" + $nbsp + "displayLabel=`"Description`">"

$ws.Range("V2").Value = $newNoteTag

# Reflect the resulting selection change recorded in the workbook (the last
# edited/selected cell moved from Y2 to V2).
$ws.Range("V2").Select()
